# Update the "function|uri list" worksheet: the Content section URLs are
# rewritten to drop the "project/<%project_id%>/" prefix (content methods
# moved out from under the project scope to avoid rewrite-parsing
# conflicts), and the explanatory note in B21 is replaced accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("function|uri list")

$ws.Range("A21").Value = "liveserver/cps/rde/rest/<%style%>/content/"

$ws.Range("A23").Value = "liveserver/cps/rde/rest/<%style%>/content/<%content_id%>/attribute/<%att_id%>"

$ws.Range("A24").Value = "liveserver/cps/rde/rest/<%style%>/content/<%content_id%>/constraints"

$ws.Range("A25").Value = "liveserver/cps/rde/rest/<%style%>/content/<%content_id%>/comments"

$ws.Range("A26").Value = "liveserver/cps/rde/rest/<%style%>/content/<%content_id%>/ratings"

$ws.Range("A27").Value = "liveserver/cps/rde/rest/<%style%>/content/<%content_id%>/hits"

$ws.Range("A22").Value = "liveserver/cps/rde/rest/<%style%>/content/?<%content_id%>"

$ws.Range("B21").Value = "Content Methods all following items are query string now due to parsing conflicts with rewriting"

# Reflect the author's final view state: scrolled back to the top of the
# sheet (no frozen/scrolled topLeftCell override) with A21 selected.
$ws.Range("A21").Select()
